$d = $word.ActiveDocument
for ($i=1; $i -le 3; $i++) {
  $c = $d.Comments(31)
  $c.Delete()
}
